# Finalize theme refactor + add style report:
# This appends the "SMALLL" company record across the linked
# Societe / Associes / Contrat sheets of the domiciliation database.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Societe sheet: clear the stale placeholder cells on row 2
#    (only the date + activities placeholder survive there), then
#    append the new company record on row 3.
# ---------------------------------------------------------------
$wsSociete = $wb.Worksheets.Item("Societe")

$wsSociete.Cells.Item(2, 1).ClearContents()   # A2
$wsSociete.Cells.Item(2, 2).ClearContents()   # B2
$wsSociete.Cells.Item(2, 3).ClearContents()   # C2
$wsSociete.Cells.Item(2, 5).ClearContents()   # E2
$wsSociete.Cells.Item(2, 6).ClearContents()   # F2
$wsSociete.Cells.Item(2, 7).ClearContents()   # G2
$wsSociete.Cells.Item(2, 8).ClearContents()   # H2

$wsSociete.Cells.Item(3, 1).Value = "SMALLL"
$wsSociete.Cells.Item(3, 2).Value = "SARL AU"
$wsSociete.Cells.Item(3, 4).Value = "20/10/2025"
$wsSociete.Cells.Item(3, 7).Value = "56  BOULEVARD MOULAY YOUSSEF 3EME ETAGE APPT 14, CASABLANCA"
$wsSociete.Cells.Item(3, 8).Value = "Casablanca"
$wsSociete.Cells.Item(3, 9).Value = "['Travaux Divers ou de Construction', 'Marchand effectuant Import Export']"

# ---------------------------------------------------------------
# 2) Associes sheet: two blank associate placeholder rows linked
#    to the new "SMALLL" company (societe_id column Q).
# ---------------------------------------------------------------
$wsAssocies = $wb.Worksheets.Item("Associes")

$wsAssocies.Cells.Item(2, 5).Value = "20/10/2025"
$wsAssocies.Cells.Item(2, 9).Value = "20/10/2025"
$wsAssocies.Cells.Item(2, 13).Value = $false
$wsAssocies.Cells.Item(2, 17).Value = "SMALLL"

$wsAssocies.Cells.Item(3, 5).Value = "20/10/2025"
$wsAssocies.Cells.Item(3, 9).Value = "20/10/2025"
$wsAssocies.Cells.Item(3, 13).Value = $false
$wsAssocies.Cells.Item(3, 17).Value = "SMALLL"

# ---------------------------------------------------------------
# 3) Contrat sheet: the blank placeholder row shifts down from
#    row 2 to row 3 (row 2 ends up empty).
# ---------------------------------------------------------------
$wsContrat = $wb.Worksheets.Item("Contrat")

$wsContrat.Cells.Item(2, 1).ClearContents()
$wsContrat.Cells.Item(2, 2).ClearContents()
$wsContrat.Cells.Item(2, 3).ClearContents()
$wsContrat.Cells.Item(2, 4).ClearContents()
$wsContrat.Cells.Item(2, 5).ClearContents()
$wsContrat.Cells.Item(2, 6).ClearContents()

$wsContrat.Cells.Item(3, 1).Value = ""
$wsContrat.Cells.Item(3, 2).Value = ""
$wsContrat.Cells.Item(3, 3).Value = ""
$wsContrat.Cells.Item(3, 4).Value = ""
$wsContrat.Cells.Item(3, 5).Value = ""
$wsContrat.Cells.Item(3, 6).Value = ""
